$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B11").Value = "1f682c4baf00039722b9d3b2a8f6431f"
$ws.Range("B34").Value = "9b5fa738b68a8c46f512c3e8ae609d3b"
$ws.Range("B44").Value = "775da89266fde57dfe7ca7c89abf5d91"
$ws.Range("B74").Value = "8a74666dc4ebb183229cedc771aa374f"
$ws.Range("B89").Value = "e5a9c26e094a5557ae9c4aa83e416d55"
$ws.Range("B99").Value = "0c473cacc596f7b80f753639d0d0ca9c"
$ws.Range("B110").Value = "8c9098805d070995ea6995c660cc73a1"
$ws.Range("B121").Value = "81667d4f5140992663fc6287a415e11f"
$ws.Range("B154").Value = "0164192226833e8b2508d9634b0ba903"
$ws.Range("B160").Value = "adf3c1215f1ec05392a34e4fcab6d818"
$ws.Range("B161").Value = "1e5c3f3bf56fea72588394470e1cc359"
$ws.Range("B162").Value = "537a5222143850acb0b8e7c2a56d1a6f"
$ws.Range("B168").Value = "bc95cae257a5ff8399d8aa38ac0096e0"
$ws.Range("B180").Value = "8e3e66726412138b9c21d57bc4009d98"
$ws.Range("B191").Value = "aec159b771e496e8cb54e48f8a239e8e"
$ws.Range("B213").Value = "f1a3da6a4991d211f4d0e18b9486ed7a"
$ws.Range("B278").Value = "9283cf6e227051ed64790cd8214746ac"
$ws.Range("B293").Value = "a7d0b31354aa502f18e0103883abbc31"
$ws.Range("B335").Value = "fa67257d9e82773e7b9d6f5b58515c14"
$ws.Range("B345").Value = "3d3502f758d76be92c0f4e2ea3201dd1"
$ws.Range("B461").Value = "060072cb4a449d58d07838c00b609f70"
$ws.Range("B480").Value = "1fd9ef0f8869fc52d6c81138b24ec41c"
$ws.Range("B506").Value = "aa1791820592e49d2dde3aff5748084a"
$ws.Range("B514").Value = "0163ad4ebad868ebcb1fb1d515410e6b"
$ws.Range("B524").Value = "b8463e643f40c14c051b7aa3e19cc647"
$ws.Range("B534").Value = "b4d216af1c0225064ccc574065e16246"
$ws.Range("B547").Value = "61c4f18193adac7d146bc75c0f680430"
$ws.Range("B553").Value = "8317bc5e1079993b6d686cc7d773b4ef"
$ws.Range("B572").Value = "5ed55f8b2ae0bd9cea467720286f267b"
$ws.Range("B584").Value = "e375d004872e7eac94fce210d9414135"
$ws.Range("B666").Value = "d0198b482e7ad0701fea272aba6657a8"
$ws.Range("B729").Value = "b4db0bd5cfe9f51ea71702c7935a8b82"
$ws.Range("B768").Value = "856d009b685edcaa25e7aebd1e4cb92c"
$ws.Range("B811").Value = "5f1e48ea2ee37ac4a0cd6534daf28e1d"
$ws.Range("B815").Value = "deeeabb02d47e448e34e5d3bbaeb8dad"
$ws.Range("B816").Value = "831b12f239db1883cfb6a62cd480eabe"
$ws.Range("B825").Value = "e0b748b7abab51601ff88878e1646e1d"
$ws.Range("B827").Value = "e72e4ad52475855fd285dd2b5bbecbd4"
$ws.Range("B874").Value = "c9c849f03081bb7a17b5eba5feebb7ea"
